$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (shifts old costs/Costs column I -> J)
$ws.Columns.Item(9).Insert()

# New column I headers: "digestion" / "Digestion"
$ws.Range("I1").Value = "digestion"
$ws.Range("I2").Value = "Digestion"

# Fill in the performance data (rows 3-8, columns E..J) which used to be #N/A errors
$data = @(
    @(0, -1, -1, -1, -1, 0),
    @(0,  0,  0,  0,  0, 0),
    @(0,  1, -1, -1, -1, 0),
    @(1,  0,  0,  0,  1, 0),
    @(0,  1, -1, -1,  0, 1),
    @(1,  0,  0,  0,  1, -1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 3 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 5 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Column widths (closest achievable values to the target 15.44140625 / 10.21875
# via the ColumnWidth property's character-width rounding)
$ws.Columns.Item(3).ColumnWidth = 14.6
$ws.Columns.Item(6).ColumnWidth = 9.3

# Selection
$ws.Range("E11").Select() | Out-Null
